$d = $word.ActiveDocument

# --- New paragraph style: "Abstract Title" (styleId AbstractTitle) ---
# Inserted right after the "Date" style / right before "Abstract" in the
# styles part; Word appends new styles at the end of the collection, which
# is immaterial to the semantic content of the style itself.
$titleStyle = $d.Styles.Add("Abstract Title", 1)          # 1 = wdStyleTypeParagraph
$titleStyle.BaseStyle = "Normal"
$titleStyle.NextParagraphStyle = "Abstract"
$titleStyle.QuickStyle = $true

$titleStyle.ParagraphFormat.KeepWithNext = $true           # <w:keepNext/>
$titleStyle.ParagraphFormat.KeepTogether = $true            # <w:keepLines/>
$titleStyle.ParagraphFormat.Alignment = 1                   # wdAlignParagraphCenter -> <w:jc w:val="center"/>
$titleStyle.ParagraphFormat.SpaceAfter = 0                  # <w:spacing w:after="0" .../>
$titleStyle.ParagraphFormat.SpaceBefore = 15                # 15pt == 300 twips -> <w:spacing .../ w:before="300"/>

$titleStyle.Font.Size = 10                                  # <w:sz w:val="20"/>   (half-points)
$titleStyle.Font.SizeBi = 10                                # <w:szCs w:val="20"/>
$titleStyle.Font.Bold = $true                                # <w:b/>

# Font.Color expects an OLE COLORREF (0xBBGGRR); build it explicitly from
# the target RRGGBB hex (345A8A) so the byte order is obviously correct.
$r = 0x34
$g = 0x5A
$b = 0x8A
$titleStyle.Font.Color = $b * 65536 + $g * 256 + $r          # <w:color w:val="345A8A"/>

# --- "Abstract" style: shrink the space-before from 300 twips to 100 ---
$abstractStyle = $d.Styles("Abstract")
$abstractStyle.ParagraphFormat.SpaceBefore = 5                # 5pt == 100 twips

# --- "ImportTok" character style: add bold + green color ---
$importTok = $d.Styles("ImportTok")
$r = 0x00
$g = 0x80
$b = 0x00
$importTok.Font.Color = $b * 65536 + $g * 256 + $r            # <w:color w:val="008000"/>
$importTok.Font.Bold = $true                                   # <w:b/>

# --- "BuiltInTok" character style: add green color ---
$builtInTok = $d.Styles("BuiltInTok")
$r = 0x00
$g = 0x80
$b = 0x00
$builtInTok.Font.Color = $b * 65536 + $g * 256 + $r            # <w:color w:val="008000"/>
